# Auto-generated edit script: updates IFRS figures for rows 2-9
# of the "company_list" sheet, per the commit "error solve ifrs list".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 8016
$ws.Range("E2").Value = 803
$ws.Range("F2").Value = 803
$ws.Range("G2").Value = 698
$ws.Range("H2").Value = 504
$ws.Range("I2").Value = 489
$ws.Range("J2").Value = 15
$ws.Range("K2").Value = 12285
$ws.Range("L2").Value = 7628
$ws.Range("M2").Value = 4657
$ws.Range("N2").Value = 4269
$ws.Range("O2").Value = 388
$ws.Range("P2").Value = 128
$ws.Range("Q2").Value = 1640
$ws.Range("R2").Value = -1316
$ws.Range("S2").Value = -8
$ws.Range("T2").Value = 1013
$ws.Range("U2").Value = 628
$ws.Range("V2").Value = 4644
$ws.Range("W2").Value = 10.02
$ws.Range("X2").Value = 6.29
$ws.Range("Y2").Value = 12.12
$ws.Range("Z2").Value = 4.48
$ws.Range("AA2").Value = 163.79
$ws.Range("AB2").Value = 3257.2
$ws.Range("AC2").Value = 1914
$ws.Range("AD2").Value = 26.09
$ws.Range("AE2").Value = 16721
$ws.Range("AF2").Value = 2.99
$ws.Range("AG2").Value = 400
$ws.Range("AH2").Value = 0.8
$ws.Range("AJ2").Value = 25533531
$ws.Range("AI2").ClearContents()

# Row 3
$ws.Range("D3").Value = 9458
$ws.Range("E3").Value = 1145
$ws.Range("F3").Value = 1145
$ws.Range("G3").Value = 1295
$ws.Range("H3").Value = 1013
$ws.Range("I3").Value = 1018
$ws.Range("J3").Value = -5
$ws.Range("K3").Value = 18448
$ws.Range("L3").Value = 12479
$ws.Range("M3").Value = 5968
$ws.Range("N3").Value = 5587
$ws.Range("O3").Value = 381
$ws.Range("P3").Value = 136
$ws.Range("Q3").Value = -1808
$ws.Range("R3").Value = -643
$ws.Range("S3").Value = 3492
$ws.Range("T3").Value = 683
$ws.Range("U3").Value = -2490
$ws.Range("V3").Value = 7963
$ws.Range("W3").Value = 12.11
$ws.Range("X3").Value = 10.71
$ws.Range("Y3").Value = 20.66
$ws.Range("Z3").Value = 6.59
$ws.Range("AA3").Value = 209.1
$ws.Range("AB3").Value = 3999.78
$ws.Range("AC3").Value = 3822
$ws.Range("AD3").Value = 11.64
$ws.Range("AE3").Value = 20488
$ws.Range("AF3").Value = 2.17
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 1.12
$ws.Range("AJ3").Value = 27270147
$ws.Range("AI3").ClearContents()

# Row 4
$ws.Range("D4").Value = 17241
$ws.Range("E4").Value = 3048
$ws.Range("F4").Value = 3048
$ws.Range("G4").Value = 2867
$ws.Range("H4").Value = 2140
$ws.Range("I4").Value = 2123
$ws.Range("J4").Value = 17
$ws.Range("K4").Value = 20744
$ws.Range("L4").Value = 11348
$ws.Range("M4").Value = 9396
$ws.Range("N4").Value = 8999
$ws.Range("O4").Value = 397
$ws.Range("P4").Value = 154
$ws.Range("Q4").Value = 1059
$ws.Range("R4").Value = -1483
$ws.Range("S4").Value = 653
$ws.Range("T4").Value = 1466
$ws.Range("U4").Value = -407
$ws.Range("V4").Value = 7429
$ws.Range("W4").Value = 17.68
$ws.Range("X4").Value = 12.41
$ws.Range("Y4").Value = 29.12
$ws.Range("Z4").Value = 10.92
$ws.Range("AA4").Value = 120.77
$ws.Range("AB4").Value = 5743.62
$ws.Range("AC4").Value = 7340
$ws.Range("AD4").Value = 6.07
$ws.Range("AE4").Value = 29294
$ws.Range("AF4").Value = 1.52
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 2.24
$ws.Range("AJ4").Value = 30718467
$ws.Range("AI4").ClearContents()

# Row 5
$ws.Range("D5").Value = 18330
$ws.Range("E5").Value = 3245
$ws.Range("F5").Value = 3245
$ws.Range("G5").Value = 2749
$ws.Range("H5").Value = 1950
$ws.Range("I5").Value = 1952
$ws.Range("J5").Value = -2
$ws.Range("K5").Value = 22722
$ws.Range("L5").Value = 11661
$ws.Range("M5").Value = 11060
$ws.Range("N5").Value = 10647
$ws.Range("O5").Value = 413
$ws.Range("P5").Value = 154
$ws.Range("Q5").Value = 1113
$ws.Range("R5").Value = -1379
$ws.Range("S5").Value = 506
$ws.Range("T5").Value = 835
$ws.Range("U5").Value = 278
$ws.Range("V5").Value = 8235
$ws.Range("W5").Value = 17.71
$ws.Range("X5").Value = 10.64
$ws.Range("Y5").Value = 19.87
$ws.Range("Z5").Value = 8.970000000000001
$ws.Range("AA5").Value = 105.44
$ws.Range("AB5").Value = 6810.79
$ws.Range("AC5").Value = 6343
$ws.Range("AD5").Value = 5.49
$ws.Range("AE5").Value = 34576
$ws.Range("AF5").Value = 1.01
$ws.Range("AG5").Value = 1100
$ws.Range("AH5").Value = 3.16
$ws.Range("AI5").Value = 17.35
$ws.Range("AJ5").Value = 30827281

# Row 6
$ws.Range("D6").Value = 17156
$ws.Range("E6").Value = 3874
$ws.Range("F6").Value = 3874
$ws.Range("G6").Value = 3377
$ws.Range("H6").Value = 2438
$ws.Range("I6").Value = 2452
$ws.Range("K6").Value = 21827
$ws.Range("L6").Value = 10749
$ws.Range("M6").Value = 11078
$ws.Range("N6").Value = 10667
$ws.Range("P6").Value = 154
$ws.Range("Q6").Value = 2561
$ws.Range("R6").Value = -1909
$ws.Range("S6").Value = -1059
$ws.Range("T6").Value = 774
$ws.Range("U6").Value = 1788
$ws.Range("V6").Value = 7466
$ws.Range("W6").Value = 22.58
$ws.Range("X6").Value = 14.21
$ws.Range("Y6").Value = 23.01
$ws.Range("Z6").Value = 10.95
$ws.Range("AA6").Value = 97.03
$ws.Range("AB6").Value = 6785.76
$ws.Range("AC6").Value = 7952
$ws.Range("AD6").Value = 3.85
$ws.Range("AE6").Value = 34612
$ws.Range("AF6").Value = 0.89
$ws.Range("AG6").Value = 1200
$ws.Range("AH6").Value = 3.92
$ws.Range("AI6").Value = 15.08
$ws.Range("AJ6").Value = 30853371

# Row 7
$ws.Range("D7").Value = 10352
$ws.Range("E7").Value = 802
$ws.Range("G7").Value = 1114
$ws.Range("H7").Value = 836
$ws.Range("I7").Value = 758
$ws.Range("K7").Value = 25982
$ws.Range("L7").Value = 13812
$ws.Range("M7").Value = 12170
$ws.Range("N7").Value = 11120
$ws.Range("P7").Value = 152
$ws.Range("Q7").Value = 956
$ws.Range("R7").Value = -2431
$ws.Range("S7").Value = 1493
$ws.Range("T7").Value = 771
$ws.Range("U7").Value = -904
$ws.Range("W7").Value = 7.75
$ws.Range("X7").Value = 8.08
$ws.Range("Y7").Value = 6.96
$ws.Range("Z7").Value = 3.5
$ws.Range("AA7").Value = 113.49
$ws.Range("AC7").Value = 2455
$ws.Range("AD7").Value = 11.67
$ws.Range("AE7").Value = 36036
$ws.Range("AF7").Value = 0.8
$ws.Range("AG7").Value = 800
$ws.Range("AH7").Value = 2.79
$ws.Range("AI7").Value = 32.6

# Row 8
$ws.Range("D8").Value = 14520
$ws.Range("E8").Value = 2098
$ws.Range("G8").Value = 1871
$ws.Range("H8").Value = 1384
$ws.Range("I8").Value = 1226
$ws.Range("K8").Value = 27636
$ws.Range("L8").Value = 14468
$ws.Range("M8").Value = 13167
$ws.Range("N8").Value = 12110
$ws.Range("P8").Value = 152
$ws.Range("Q8").Value = 2435
$ws.Range("R8").Value = -1592
$ws.Range("S8").Value = 116
$ws.Range("T8").Value = 977
$ws.Range("U8").Value = 917
$ws.Range("W8").Value = 14.45
$ws.Range("X8").Value = 9.529999999999999
$ws.Range("Y8").Value = 10.55
$ws.Range("Z8").Value = 5.16
$ws.Range("AA8").Value = 109.88
$ws.Range("AC8").Value = 3967
$ws.Range("AD8").Value = 7.22
$ws.Range("AE8").Value = 39245
$ws.Range("AF8").Value = 0.73
$ws.Range("AG8").Value = 1200
$ws.Range("AH8").Value = 4.19
$ws.Range("AI8").Value = 30.25

# Row 9
$ws.Range("D9").Value = 16918
$ws.Range("E9").Value = 2939
$ws.Range("G9").Value = 2754
$ws.Range("H9").Value = 2034
$ws.Range("I9").Value = 1848
$ws.Range("K9").Value = 30000
$ws.Range("L9").Value = 15320
$ws.Range("M9").Value = 14680
$ws.Range("N9").Value = 13750
$ws.Range("P9").Value = 152
$ws.Range("Q9").Value = -141
$ws.Range("R9").Value = -1404
$ws.Range("S9").Value = -256
$ws.Range("T9").Value = 1195
$ws.Range("U9").Value = -1498
$ws.Range("W9").Value = 17.37
$ws.Range("X9").Value = 12.03
$ws.Range("Y9").Value = 14.3
$ws.Range("Z9").Value = 7.06
$ws.Range("AA9").Value = 104.36
$ws.Range("AC9").Value = 5984
$ws.Range("AD9").Value = 4.79
$ws.Range("AE9").Value = 44559
$ws.Range("AF9").Value = 0.64
$ws.Range("AG9").Value = 1450
$ws.Range("AH9").Value = 5.06
$ws.Range("AI9").Value = 24.23

